# Appended/Updated sheet 'HOTCARD' via Github_Excel_Append_Tool
#
# - Adds 9 new trailing header columns (CU1:DC1) on the HOTCARD sheet.
# - Retypes several existing data cells in rows 2-3 from text to numeric
#   (ZIP code, routing #, BIN, BIN length, port, card count columns).
# - Leaves the new trailing columns blank for the two existing data rows.
# - Appends a new row 4 with an Entity ID (K4) and values for the new
#   trailing columns (CU4:DC4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOTCARD")

# ---------------------------------------------------------------------
# 1) New header cells CU1:DC1 (bold/bordered header style, same as the
#    rest of row 1 - copy format from the existing CT1 header cell).
# ---------------------------------------------------------------------
$headers = @{
    "CU1" = "Migration Date"
    "CV1" = "FI Name"
    "CW1" = "Switch: FISB"
    "CX1" = "Old Platform"
    "CY1" = "New Platform :PaymentsOne Debit"
    "CZ1" = "Service: Basic"
    "DA1" = "CS Location: Offshore"
    "DB1" = "Total Card Count"
    "DC1" = "Using OneCall IVR"
}

$ct1 = $ws.Range("CT1")
foreach ($addr in $headers.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $headers[$addr]
    $ct1.Copy()
    $cell.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 2) Retype existing inline-string numeric-looking values as real
#    numbers in rows 2 and 3.
# ---------------------------------------------------------------------
$numericCells = "R2","AP2","AY2","AZ2","BP2","CR2","R3","AP3","AY3","AZ3","BP3","CR3"
$numericValues = @{
    "R2" = 62701;     "R3" = 62701
    "AP2" = 123456789; "AP3" = 123456789
    "AY2" = 123456;    "AY3" = 123456
    "AZ2" = 6;         "AZ3" = 6
    "BP2" = 8080;      "BP3" = 8080
    "CR2" = 15000;     "CR3" = 15000
}
foreach ($addr in $numericCells) {
    $ws.Range($addr).Value = $numericValues[$addr]
}

# ---------------------------------------------------------------------
# 3) New row 4: Entity ID in K4, plus values for the new trailing
#    columns CU4:DC4. The remaining cells on row 4 stay blank, matching
#    the blank placeholder cells left across A4:CT4 (excluding K4).
# ---------------------------------------------------------------------
$ws.Range("K4").Value = "123abx007"

$row4 = @{
    "CU4" = "NA"
    "CV4" = "NA"
    "CW4" = "FISB"
    "CX4" = "NA"
    "CY4" = "NA"
    "CZ4" = "NA"
    "DA4" = "NA"
    "DB4" = "NA"
    "DC4" = "NA"
}
foreach ($addr in $row4.Keys) {
    $ws.Range($addr).Value = $row4[$addr]
}
